$wb = $excel.ActiveWorkbook

# --- Sheet: BC por zonas (B/C columns, rows 2-20) ---
$ws7 = $wb.Worksheets.Item("BC por zonas")
$ws7.Range("B2").Value = 3426
$ws7.Range("C2").Value = 5193
$ws7.Range("B3").Value = 1140
$ws7.Range("C3").Value = 184
$ws7.Range("B4").Value = 1499
$ws7.Range("C4").Value = 467
$ws7.Range("B5").Value = 329
$ws7.Range("C5").Value = 11
$ws7.Range("B6").Value = 1625
$ws7.Range("C6").Value = 2378
$ws7.Range("B7").Value = 1649
$ws7.Range("C7").Value = 2661
$ws7.Range("B8").Value = 182
$ws7.Range("C8").Value = 149
$ws7.Range("B9").Value = 326
$ws7.Range("C9").Value = 135
$ws7.Range("B10").Value = 103
$ws7.Range("C10").Value = 30
$ws7.Range("B11").Value = 1171
$ws7.Range("C11").Value = 1072
$ws7.Range("B12").Value = 1225
$ws7.Range("C12").Value = 3162
$ws7.Range("B13").Value = 230
$ws7.Range("C13").Value = 136
$ws7.Range("B14").Value = 120
$ws7.Range("C14").Value = 305
$ws7.Range("B15").Value = 596
$ws7.Range("C15").Value = 333
$ws7.Range("B16").Value = 601
$ws7.Range("C16").Value = 243
$ws7.Range("B17").Value = 649
$ws7.Range("C17").Value = 102
$ws7.Range("B18").Value = 68
$ws7.Range("C18").Value = 52
$ws7.Range("B19").Value = 132
$ws7.Range("C19").Value = 91
$ws7.Range("B20").Value = 814
$ws7.Range("C20").Value = 474
[void]$ws7.Range("F16").Select()

# --- Sheet: Expo-ICA (B column) ---
$ws8 = $wb.Worksheets.Item("Expo-ICA")
$ws8.Range("B2").Value = 15888
$ws8.Range("B3").Value = 3651
$ws8.Range("B4").Value = 7
$ws8.Range("B5").Value = 408
$ws8.Range("B6").Value = 46
$ws8.Range("B7").Value = 158
$ws8.Range("B8").Value = 98
$ws8.Range("B9").Value = 2374
$ws8.Range("B10").Value = 360
$ws8.Range("B11").Value = 63
$ws8.Range("B12").Value = 6
$ws8.Range("B13").Value = 2
$ws8.Range("B14").Value = 86
$ws8.Range("B15").Value = 42
$ws8.Range("B16").Value = 5514
$ws8.Range("B17").Value = 802
$ws8.Range("B18").Value = 60
$ws8.Range("B19").Value = 323
$ws8.Range("B20").Value = 11
$ws8.Range("B21").Value = 15
$ws8.Range("B22").Value = 37
$ws8.Range("B23").Value = 229
$ws8.Range("B24").Value = 1428
$ws8.Range("B25").Value = 51
$ws8.Range("B26").Value = 168
$ws8.Range("B27").Value = 181
$ws8.Range("B28").Value = 1916
$ws8.Range("B29").Value = 60
$ws8.Range("B30").Value = 96
$ws8.Range("B31").Value = 23
$ws8.Range("B32").Value = 111
$ws8.Range("B33").Value = 4719
$ws8.Range("B34").Value = 933
$ws8.Range("B35").Value = 203
$ws8.Range("B36").Value = 40
$ws8.Range("B37").Value = 5
$ws8.Range("B38").Value = 95
$ws8.Range("B39").Value = 24
$ws8.Range("B41").Value = 41
$ws8.Range("B42").Value = 660
$ws8.Range("B43").Value = 519
$ws8.Range("B44").Value = 365
$ws8.Range("B45").Value = 1715
$ws8.Range("B46").Value = 52
$ws8.Range("B47").Value = 65
$ws8.Range("B48").Value = 2003
$ws8.Range("B49").Value = 786
$ws8.Range("B50").Value = 602
$ws8.Range("B51").Value = 17
$ws8.Range("B52").Value = 472
$ws8.Range("B53").Value = 127

# --- Sheet: Impo-ICA (B column) ---
$ws9 = $wb.Worksheets.Item("Impo-ICA")
$ws9.Range("B2").Value = 17178
$ws9.Range("B3").Value = 2578
$ws9.Range("B4").Value = 1909
$ws9.Range("B5").Value = 261
$ws9.Range("B6").Value = 322
$ws9.Range("B7").Value = 87
$ws9.Range("B8").Value = 6673
$ws9.Range("B9").Value = 1184
$ws9.Range("B10").Value = 293
$ws9.Range("B11").Value = 2047
$ws9.Range("B12").Value = 842
$ws9.Range("B13").Value = 277
$ws9.Range("B14").Value = 257
$ws9.Range("B15").Value = 128
$ws9.Range("B16").Value = 1198
$ws9.Range("B17").Value = 446
$ws9.Range("B18").Value = 1810
$ws9.Range("B19").Value = 1706
$ws9.Range("B20").Value = 104
$ws9.Range("B21").Value = 3640
$ws9.Range("B22").Value = 2296
$ws9.Range("B23").Value = 1102
$ws9.Range("B24").Value = 242
$ws9.Range("B25").Value = 1957
$ws9.Range("B26").Value = 47
$ws9.Range("B27").Value = 158
$ws9.Range("B28").Value = 191
$ws9.Range("B29").Value = 621
$ws9.Range("B30").Value = 89
$ws9.Range("B31").Value = 134
$ws9.Range("B32").Value = 176
$ws9.Range("B33").Value = 112
$ws9.Range("B34").Value = 125
$ws9.Range("B35").Value = 73
$ws9.Range("B36").Value = 106
$ws9.Range("B37").Value = 127
$ws9.Range("B38").Value = 414
$ws9.Range("B39").Value = 414
$ws9.Range("B40").Value = 106
[void]$ws9.Range("A1").Select()

# --- Activate Expo-ICA as the final active sheet/tab ---
[void]$ws8.Select()
